# Add BMC ports to river compute and add tests.
#
# HARDWARE_MANAGEMENT (sheet index 3) gets 4 new rows (25-28) describing
# the BMC port connections for cn01-cn04.
# COMPUTE_NODES (sheet index 4) gets the matching 4 new rows (24-27).
# The active tab moves from INTER_SWITCH_LINKS (sheet 1) to COMPUTE_NODES
# (sheet 4), and the selections on HARDWARE_MANAGEMENT / COMPUTE_NODES are
# updated to the newly added ranges.

$wb = $excel.ActiveWorkbook

$wsHardwareManagement = $wb.Worksheets.Item(3)
$wsComputeNodes = $wb.Worksheets.Item(4)

# --- HARDWARE_MANAGEMENT: new rows 25-28 (BMC ports for cn01..cn04) ---
$bmcRows = @(
    @{ Row = 25; Source = "cn01"; Location = "u15"; Port = 11 },
    @{ Row = 26; Source = "cn02"; Location = "u16"; Port = 12 },
    @{ Row = 27; Source = "cn03"; Location = "u17"; Port = 13 },
    @{ Row = 28; Source = "cn04"; Location = "u18"; Port = 14 }
)

foreach ($r in $bmcRows) {
    $row = $r.Row
    $wsHardwareManagement.Range("J$row").Value = $r.Source
    $wsHardwareManagement.Range("K$row").Value = "x3002"
    $wsHardwareManagement.Range("L$row").Value = $r.Location
    $wsHardwareManagement.Range("M$row").Value = "bmc"
    $wsHardwareManagement.Range("O$row").Value = 1
    $wsHardwareManagement.Range("P$row").Value = "sw-leaf-bmc-001"
    $wsHardwareManagement.Range("Q$row").Value = "x3000"
    $wsHardwareManagement.Range("R$row").Value = "u37"
    $wsHardwareManagement.Range("T$row").Value = $r.Port
}

# --- COMPUTE_NODES: new rows 24-27 (BMC ports for cn01..cn04) ---
$computeRows = @(
    @{ Row = 24; Source = "cn01"; Location = "u15"; Port = 24 },
    @{ Row = 25; Source = "cn02"; Location = "u16"; Port = 25 },
    @{ Row = 26; Source = "cn03"; Location = "u17"; Port = 26 },
    @{ Row = 27; Source = "cn04"; Location = "u18"; Port = 27 }
)

foreach ($r in $computeRows) {
    $row = $r.Row
    $wsComputeNodes.Range("J$row").Value = $r.Source
    $wsComputeNodes.Range("K$row").Value = "x3002"
    $wsComputeNodes.Range("L$row").Value = $r.Location
    $wsComputeNodes.Range("O$row").Value = 1
    $wsComputeNodes.Range("P$row").Value = "sw-leaf-bmc-001"
    $wsComputeNodes.Range("Q$row").Value = "x3000"
    $wsComputeNodes.Range("R$row").Value = "u37"
    $wsComputeNodes.Range("T$row").Value = $r.Port
}

# --- Update selections to reflect the newly added rows ---
[void]$wsHardwareManagement.Activate()
[void]$wsHardwareManagement.Range("J25:T28").Select()

# --- Move the active/selected tab to COMPUTE_NODES ---
[void]$wsComputeNodes.Activate()
[void]$wsComputeNodes.Range("J24:T27").Select()
